$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: fill in the new item's data (counter, name, ratio, qty, price, sell price, transactions)
$ws.Range("A7").Value = 1

$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "CONCOR PLUS 10/25MG 30 F.C. TABLETS"

$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "1:2"

$ws.Range("L7:M7").NumberFormat = "@"
$ws.Range("L7").Value = "1"

$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "108.00"

$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "35.6400"

$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "0:1"

# Row 8: total line, now shows the computed total and matches row 7's height
$ws.Range("N8").Value = 35.64
$ws.Rows.Item(8).RowHeight = 25.5
